$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold "price" text which can look like a plain number (e.g. "311.12")
# are force-formatted as Text before the write so Excel keeps them as strings
# (matching the source t="inlineStr" cells) instead of silently coercing them
# into numeric values; ClearFormats() afterwards drops the explicit format again
# so the cell style stays at the sheet default, same as every other data cell.
$textCells = @{
    'D2' = '43.369.49'
    'D3' = '2.310.29'
    'D5' = '311.12'
    'D6' = '102.98'
    'D10' = '35.73'
    'D12' = '52.09'
    'D14' = '7.03'
    'D15' = '2.671.02'
    'D16' = '15.01'
    'D17' = '2.314.43'
    'D19' = '43.283.16'
    'D20' = '12.27'
    'D21' = '0.0₃0932'
    'D23' = '68.12'
    'D24' = '241.65'
    'D28' = '24.89'
    'D29' = '2.30'
    'D30' = '36.80'
    'D31' = '9.65'
    'D32' = '169.92'
    'D36' = '17.71'
    'D37' = '0.0742'
    'D43' = '2.38'
    'D45' = '1.977.05'
    'D46' = '19.16'
    'D48' = '10.00'
    'D49' = '2.95'
    'D50' = '55.41'
    'D51' = '1.60'
}
foreach ($ref in $textCells.Keys) {
    $range = $ws.Range($ref)
    $range.NumberFormat = '@'
    $range.Value = $textCells[$ref]
    $range.ClearFormats()
}

# Remaining plain-text cell updates (coin names, links, percentage strings).
$ws.Range('E2').Value = '  +3.17%  '
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +1.68%  '
$ws.Range('E6').Value = '  +6.46%  '
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +8.66%  '
$ws.Range('E10').Value = '  +2.49%  '
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('E16').Value = '  +2.80%  '
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('E18').Value = '  +2.97%  '
$ws.Range('E19').Value = '  +3.34%  '
$ws.Range('E20').Value = '  +0.68%  '
$ws.Range('E21').Value = '  +3.73%  '
$ws.Range('E22').Value = '  +3.51%  '
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('E25').Value = '  +2.75%  '
$ws.Range('E26').Value = '  +2.05%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +6.34%  '
$ws.Range('E29').Value = '  +8.83%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('E32').Value = '  +3.49%  '
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +6.82%  '
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('E42').Value = '  +5.94%  '
$ws.Range('E43').Value = '  +1.98%  '
$ws.Range('E44').Value = '  +4.80%  '
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('E47').Value = '  +4.22%  '
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E49').Value = '  +2.82%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('E51').Value = '  +10.31%  '
